{"js": "const replacements = [\n  [\"93\u00d750=4650\", \"14\u00d739=546\"],\n  [\"34\u00d712=408\", \"85\u00d756=4760\"],\n  [\"77\u00d774=5698\", \"21\u00d774=1554\"],\n  [\"94\u00d766=6204\", \"99\u00d792=9108\"],\n  [\"29\u00d777=2233\", \"44\u00d781=3564\"],\n  [\"98\u00d770=6860\", \"40\u00d760=2400\"],\n  [\"38\u00d762=2356\", \"59\u00d758=3422\"],\n  [\"59\u00d716=944\", \"29\u00d728=812\"],\n  [\"70\u00d765=4550\", \"87\u00d744=3828\"],\n  [\"41\u00d771=2911\", \"84\u00d743=3612\"],\n  [\"11\u00d746=506\", \"47\u00d742=1974\"],\n  [\"26\u00d733=858\", \"77\u00d788=6776\"],\n  [\"30\u00d724=720\", \"43\u00d724=1032\"],\n  [\"64\u00d717=1088\", \"19\u00d722=418\"],\n  [\"75\u00d760=4500\", \"66\u00d782=5412\"],\n  [\"99\u00d738=3762\", \"21\u00d769=1449\"],\n  [\"16\u00d775=1200\", \"86\u00d770=6020\"],\n  [\"44\u00d731=1364\", \"15\u00d782=1230\"],\n  [\"37\u00d717=629\", \"15\u00d781=1215\"],\n  [\"24\u00d728=672\", \"86\u00d794=8084\"],\n  [\"87\u00d720=1740\", \"11\u00d785=935\"],\n  [\"47\u00d797=4559\", \"35\u00d713=455\"],\n  [\"36\u00d727=972\", \"66\u00d742=2772\"],\n  [\"58\u00d771=4118\", \"66\u00d716=1056\"],\n  [\"33\u00d746=1518\", \"56\u00d791=5096\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  // Replace every match (expected to be exactly one occurrence each).\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"93\u00d750=4650\"; New = \"14\u00d739=546\" },\n    @{ Old = \"34\u00d712=408\"; New = \"85\u00d756=4760\" },\n    @{ Old = \"77\u00d774=5698\"; New = \"21\u00d774=1554\" },\n    @{ Old = \"94\u00d766=6204\"; New = \"99\u00d792=9108\" },\n    @{ Old = \"29\u00d777=2233\"; New = \"44\u00d781=3564\" },\n    @{ Old = \"98\u00d770=6860\"; New = \"40\u00d760=2400\" },\n    @{ Old = \"38\u00d762=2356\"; New = \"59\u00d758=3422\" },\n    @{ Old = \"59\u00d716=944\"; New = \"29\u00d728=812\" },\n    @{ Old = \"70\u00d765=4550\"; New = \"87\u00d744=3828\" },\n    @{ Old = \"41\u00d771=2911\"; New = \"84\u00d743=3612\" },\n    @{ Old = \"11\u00d746=506\"; New = \"47\u00d742=1974\" },\n    @{ Old = \"26\u00d733=858\"; New = \"77\u00d788=6776\" },\n    @{ Old = \"30\u00d724=720\"; New = \"43\u00d724=1032\" },\n    @{ Old = \"64\u00d717=1088\"; New = \"19\u00d722=418\" },\n    @{ Old = \"75\u00d760=4500\"; New = \"66\u00d782=5412\" },\n    @{ Old = \"99\u00d738=3762\"; New = \"21\u00d769=1449\" },\n    @{ Old = \"16\u00d775=1200\"; New = \"86\u00d770=6020\" },\n    @{ Old = \"44\u00d731=1364\"; New = \"15\u00d782=1230\" },\n    @{ Old = \"37\u00d717=629\"; New = \"15\u00d781=1215\" },\n    @{ Old = \"24\u00d728=672\"; New = \"86\u00d794=8084\" },\n    @{ Old = \"87\u00d720=1740\"; New = \"11\u00d785=935\" },\n    @{ Old = \"47\u00d797=4559\"; New = \"35\u00d713=455\" },\n    @{ Old = \"36\u00d727=972\"; New = \"66\u00d742=2772\" },\n    @{ Old = \"58\u00d771=4118\"; New = \"66\u00d716=1056\" },\n    @{ Old = \"33\u00d746=1518\"; New = \"56\u00d791=5096\" },\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
